$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "04-21-2024"

# Header row
$ws.Range("B1").Value = "Linear Predicted Difference"
$ws.Range("C1").Value = "XGBoost Predicted Difference"

# Data rows (A2:C37)
$data = New-Object 'object[,]' 36,3
$data[0,0] = "pollster_CBS News"
$data[0,1] = 3.072917012648367
$data[0,2] = 2.549972295761108
$data[1,0] = "pollster_CNBC"
$data[1,1] = 1.142966447193374
$data[1,2] = 0.8700668811798096
$data[2,0] = "pollster_CNN"
$data[2,1] = 3.686799599260403
$data[2,2] = 2.971726894378662
$data[3,0] = "pollster_Daily Kos/Civiqs"
$data[3,1] = 0.182212138579338
$data[3,2] = 0.9383037686347961
$data[4,0] = "pollster_Data for Progress (D)**"
$data[4,1] = -1.207759762355517
$data[4,2] = -0.3364807069301605
$data[5,0] = "pollster_Economist/YouGov"
$data[5,1] = 0.6908655814551774
$data[5,2] = 0.7516812086105347
$data[6,0] = "pollster_Emerson"
$data[6,1] = 2.159319310666717
$data[6,2] = 2.989696025848389
$data[7,0] = "pollster_FOX News"
$data[7,1] = 2.484873686204431
$data[7,2] = 4.331109046936035
$data[8,0] = "pollster_Federalist/Susquehanna"
$data[8,1] = -1.608179587621299
$data[8,2] = -2.487954378128052
$data[9,0] = "pollster_Forbes/HarrisX"
$data[9,1] = 2.841870646069874
$data[9,2] = 0.3763706684112549
$data[10,0] = "pollster_Grinnell/Selzer"
$data[10,1] = 2.913608913116613
$data[10,2] = 6.060928344726562
$data[11,0] = "pollster_HarrisX**"
$data[11,1] = 6.264908416857561
$data[11,2] = 2.607294082641602
$data[12,0] = "pollster_Harvard-Harris"
$data[12,1] = 5.355491091275961
$data[12,2] = 2.15910530090332
$data[13,0] = "pollster_I&I/TIPP"
$data[13,1] = 0.3238862582728297
$data[13,2] = -2.186527252197266
$data[14,0] = "pollster_Insider Advantage"
$data[14,1] = 5.534786859571602
$data[14,2] = 4.231391429901123
$data[15,0] = "pollster_InsiderAdvantage"
$data[15,1] = 4.319954623606749
$data[15,2] = 7.111198425292969
$data[16,0] = "pollster_Marist"
$data[16,1] = 0.05963366743036502
$data[16,2] = 0.717138946056366
$data[17,0] = "pollster_Marquette"
$data[17,1] = 2.774417667570467
$data[17,2] = 0.8810979127883911
$data[18,0] = "pollster_Morning Consult"
$data[18,1] = 0.7243218602052033
$data[18,2] = 0.01597485318779945
$data[19,0] = "pollster_NBC News"
$data[19,1] = 1.649235866076296
$data[19,2] = 1.948659896850586
$data[20,0] = "pollster_NPR/PBS/Marist"
$data[20,1] = -1.258409458868962
$data[20,2] = -1.311262249946594
$data[21,0] = "pollster_NY Times/Siena"
$data[21,1] = 0.4165032637335209
$data[21,2] = 1.131097793579102
$data[22,0] = "pollster_New York Post"
$data[22,1] = 3.943999041024175
$data[22,2] = 3.40349555015564
$data[23,0] = "pollster_PPP (D)"
$data[23,1] = -2.166665559525594
$data[23,2] = -1.257511377334595
$data[24,0] = "pollster_Politico/Morning Consult"
$data[24,1] = 0.7198022618213797
$data[24,2] = 0.2121571153402328
$data[25,0] = "pollster_Quinnipiac"
$data[25,1] = -1.444926229922307
$data[25,2] = -2.618772268295288
$data[26,0] = "pollster_Rasmussen Reports"
$data[26,1] = 6.431448004655588
$data[26,2] = 7.59668493270874
$data[27,0] = "pollster_Reuters/Ipsos"
$data[27,1] = 0.7854871219188557
$data[27,2] = -3.087626457214355
$data[28,0] = "pollster_SurveyUSA"
$data[28,1] = 0.8493888921466355
$data[28,2] = 0.8449379205703735
$data[29,0] = "pollster_Susquehanna"
$data[29,1] = -11.22973899211079
$data[29,2] = -7.595902919769287
$data[30,0] = "pollster_The Messenger/HarrisX"
$data[30,1] = 4.104119399590374
$data[30,2] = 4.241454601287842
$data[31,0] = "pollster_Trafalgar Group (R)"
$data[31,1] = 6.189985428106216
$data[31,2] = 5.265925884246826
$data[32,0] = "pollster_USA Today/Suffolk"
$data[32,1] = -0.5415133638111804
$data[32,2] = 0.6435413956642151
$data[33,0] = "pollster_Wall Street Journal"
$data[33,1] = 0.2313121291008491
$data[33,2] = 0.6309971809387207
$data[34,0] = "pollster_Yahoo News"
$data[34,1] = -0.302828325158921
$data[34,2] = 0.1194566562771797
$data[35,0] = "pollster_Yahoo News**"
$data[35,1] = -1.714591887196505
$data[35,2] = -3.566738605499268
$ws.Range("A2:C37").Value = $data

# Styling for header row (B1:C1) and pollster name column (A2:A37): bold font, thin border, center/top alignment
$styledRange = $ws.Range("B1:C1")
$styledRange.Font.Bold = $true
$styledRange.HorizontalAlignment = -4108
$styledRange.VerticalAlignment = -4160
$styledRange.Borders.LineStyle = 1

$styledRange2 = $ws.Range("A2:A37")
$styledRange2.Font.Bold = $true
$styledRange2.HorizontalAlignment = -4108
$styledRange2.VerticalAlignment = -4160
$styledRange2.Borders.LineStyle = 1

# Page margins to match: left/right 0.75in, top/bottom 1in, header/footer 0.5in
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
